$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Header label text shortened
Replace-Text "Atlantic cod" "cod"
Replace-Text "European hake" "hake"

# Numeric values rounded to 3 decimal places
Replace-Text "0.9589" "0.959"
Replace-Text "0.9543" "0.954"
Replace-Text "0.9078" "0.908"
Replace-Text "0.8973" "0.897"
Replace-Text "0.8398" "0.840"
Replace-Text "0.9889" "0.989"
Replace-Text "0.7340" "0.734"
Replace-Text "0.9043" "0.904"
Replace-Text "0.6039" "0.604"
Replace-Text "0.9392" "0.939"
Replace-Text "0.8266" "0.827"
Replace-Text "0.6965" "0.697"
Replace-Text "0.6923" "0.692"
Replace-Text "0.9958" "0.996"
Replace-Text "0.5539" "0.554"
Replace-Text "0.6497" "0.650"
Replace-Text "0.7732" "0.773"
Replace-Text "0.6869" "0.687"
Replace-Text "0.6466" "0.647"
Replace-Text "0.5666" "0.567"
